# Rename the two worksheets and update the active selection on the
# second sheet, per the commit "cleaned up test data".

$wb = $excel.ActiveWorkbook

# Rename worksheets (order/rId unaffected, only the tab names change)
$wsTwo = $wb.Worksheets.Item("test_new_survey_import 2")
$wsTwo.Name = "new_survey_import_2_test"

$wsThree = $wb.Worksheets.Item("test_new_survey_import 3")
$wsThree.Name = "new_survey_import_3_test"

# Update the saved selection/active cell on the "3" sheet (it is the
# active/selected tab in the workbook) from F18 to E11
$wsThree.Activate()
$wsThree.Range("E11").Select()
